$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I106").Value = 1751.25
$ws.Range("J106").Value = 1850
$ws.Range("K106").Value = 1751.25
$ws.Range("L106").Value = 1850
$ws.Range("M106").Value = -1120.25
$ws.Range("N106").Value = -3112

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3406.2974
$ws.Range("I74").Value = 702.84375
$ws.Range("J74").Value = 20708.4
$ws.Range("K74").Value = 702.84375
$ws.Range("L74").Value = 20708.4
$ws.Range("M74").Value = 171.15625
$ws.Range("N74").Value = -22456.4
$ws.Range("H77").Value = 3406.2974
$ws.Range("I77").Value = 702.84375
$ws.Range("J77").Value = 20708.4
$ws.Range("K77").Value = 3514.21875
$ws.Range("L77").Value = 103542
$ws.Range("M77").Value = 853.78125
$ws.Range("N77").Value = -112278

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1332.075
$ws.Range("I107").Value = 980.55
$ws.Range("J107").Value = 1683.6
$ws.Range("K107").Value = 980.55
$ws.Range("L107").Value = 1683.6
$ws.Range("M107").Value = 939.45
$ws.Range("N107").Value = -5523.6

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1660.0667
$ws.Range("I16").Value = 900.125
$ws.Range("J16").Value = 2528.5715
$ws.Range("K16").Value = 900.125
$ws.Range("L16").Value = 2528.5715
$ws.Range("M16").Value = -613.125
$ws.Range("N16").Value = -3102.5715
$ws.Range("H86").Value = 3886.5625
$ws.Range("I86").Value = 3703.4
$ws.Range("J86").Value = 3969.818
$ws.Range("K86").Value = 3703.4
$ws.Range("L86").Value = 3969.818
$ws.Range("M86").Value = -2580.4
$ws.Range("N86").Value = -6215.818
$ws.Range("H89").Value = 3886.5625
$ws.Range("I89").Value = 3703.4
$ws.Range("J89").Value = 3969.818
$ws.Range("K89").Value = 18517
$ws.Range("L89").Value = 19849.09
$ws.Range("M89").Value = -12901
$ws.Range("N89").Value = -31081.09
$ws.Range("H113").Value = 1660.0667
$ws.Range("I113").Value = 900.125
$ws.Range("J113").Value = 2528.5715
$ws.Range("K113").Value = 900.125
$ws.Range("L113").Value = 2528.5715
$ws.Range("M113").Value = 1269.875
$ws.Range("N113").Value = -6868.5715

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H58").Value = 166667180
$ws.Range("J58").Value = 166667180
$ws.Range("L58").Value = 500001540
$ws.Range("N58").Value = -500001796
$ws.Range("H70").Value = 5538
$ws.Range("I70").Value = 945
$ws.Range("J70").Value = 8600
$ws.Range("K70").Value = 2835
$ws.Range("L70").Value = 25800
$ws.Range("M70").Value = -2520
$ws.Range("N70").Value = -26430
$ws.Range("H73").Value = 5538
$ws.Range("I73").Value = 945
$ws.Range("J73").Value = 8600
$ws.Range("K73").Value = 2835
$ws.Range("L73").Value = 25800
$ws.Range("M73").Value = -1743
$ws.Range("N73").Value = -27984
$ws.Range("H76").Value = 650
$ws.Range("I76").Value = 650
$ws.Range("J76").Value = 0
$ws.Range("K76").Value = 1950
$ws.Range("L76").Value = 0
$ws.Range("M76").Value = -1567
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 650
$ws.Range("I79").Value = 650
$ws.Range("J79").Value = 0
$ws.Range("K79").Value = 1950
$ws.Range("L79").Value = 0
$ws.Range("M79").Value = -624
$ws.Range("N79").ClearContents()
$ws.Range("H107").Value = 1033.125
$ws.Range("I107").Value = 257.66666
$ws.Range("J107").Value = 1498.4
$ws.Range("K107").Value = 772.9999799999999
$ws.Range("L107").Value = 4495.200000000001
$ws.Range("M107").Value = 1147.00002
$ws.Range("N107").Value = -8335.200000000001
$ws.Range("H114").Value = 549.2083
$ws.Range("I114").Value = 208.5
$ws.Range("J114").Value = 719.5625
$ws.Range("K114").Value = 625.5
$ws.Range("L114").Value = 2158.6875
$ws.Range("M114").Value = 2628.5
$ws.Range("N114").Value = -8666.6875
$ws.Range("H121").Value = 5664.25
$ws.Range("I121").Value = 432.85715
$ws.Range("J121").Value = 6653.973
$ws.Range("K121").Value = 1298.57145
$ws.Range("L121").Value = 19961.919
$ws.Range("M121").Value = 11.42855000000009
$ws.Range("N121").Value = -22581.919
$ws.Range("H129").Value = 958.25
$ws.Range("I129").Value = 700
$ws.Range("J129").Value = 1216.5
$ws.Range("K129").Value = 2100
$ws.Range("L129").Value = 3649.5
$ws.Range("M129").Value = 2900
$ws.Range("N129").Value = -13649.5

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H38").Value = 10021
$ws.Range("I38").Value = 10021
$ws.Range("K38").Value = 10021
$ws.Range("M38").Value = -9558
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H43").Value = 13125
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 13125
$ws.Range("K43").Value = 0
$ws.Range("L43").Value = 13125
$ws.Range("M43").ClearContents()
$ws.Range("N43").Value = -13427
$ws.Range("H44").Value = 0
$ws.Range("I44").Value = 0
$ws.Range("K44").Value = 0
$ws.Range("M44").ClearContents()
$ws.Range("H46").Value = 7150
$ws.Range("I46").Value = 850
$ws.Range("J46").Value = 13450
$ws.Range("K46").Value = 850
$ws.Range("L46").Value = 13450
$ws.Range("M46").Value = -694
$ws.Range("N46").Value = -13762
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H49").Value = 6000
$ws.Range("J49").Value = 6000
$ws.Range("L49").Value = 6000
$ws.Range("N49").Value = -6368
$ws.Range("H80").Value = 3029.2666
$ws.Range("I80").Value = 2390
$ws.Range("J80").Value = 3588.625
$ws.Range("K80").Value = 2390
$ws.Range("L80").Value = 3588.625
$ws.Range("M80").Value = -1392
$ws.Range("N80").Value = -5584.625
$ws.Range("H83").Value = 3029.2666
$ws.Range("I83").Value = 2390
$ws.Range("J83").Value = 3588.625
$ws.Range("K83").Value = 11950
$ws.Range("L83").Value = 17943.125
$ws.Range("M83").Value = -6958
$ws.Range("N83").Value = -27927.125
$ws.Range("H113").Value = 41038.72
$ws.Range("I113").Value = 67407.47
$ws.Range("J113").Value = 1485.6
$ws.Range("K113").Value = 67407.47
$ws.Range("L113").Value = 1485.6
$ws.Range("M113").Value = -65237.47
$ws.Range("N113").Value = -5825.6

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 172.77272
$ws.Range("I55").Value = 153
$ws.Range("J55").Value = 189.25
$ws.Range("K55").Value = 153
$ws.Range("L55").Value = 189.25
$ws.Range("M55").Value = 20
$ws.Range("N55").Value = -535.25
$ws.Range("H64").Value = 36574.75
$ws.Range("J64").Value = 36574.75
$ws.Range("L64").Value = 36574.75
$ws.Range("N64").Value = -37024.75
$ws.Range("H67").Value = 36574.75
$ws.Range("J67").Value = 36574.75
$ws.Range("L67").Value = 36574.75
$ws.Range("N67").Value = -38134.75

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 36624.5
$ws.Range("J63").Value = 36624.5
$ws.Range("L63").Value = 36624.5
$ws.Range("N63").Value = -37872.5
$ws.Range("H66").Value = 36624.5
$ws.Range("J66").Value = 36624.5
$ws.Range("L66").Value = 109873.5
$ws.Range("N66").Value = -116113.5
$ws.Range("H113").Value = 80273.24000000001
$ws.Range("I113").Value = 58966.293
$ws.Range("J113").Value = 125550.5
$ws.Range("K113").Value = 176898.879
$ws.Range("L113").Value = 376651.5
$ws.Range("M113").Value = -174728.879
$ws.Range("N113").Value = -380991.5
